# removed hydro pump and changed it to 12h bess in xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "battery_12h_mw"
$ws.Range("R1").Value = "battery_12h_price"

$ws.Range("R2").Select()
